$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 32
$ws.Range("D32").Value = "PCA (Principal Component Analysis) 종류"
$ws.Range("E32").Value = "https://dodonam.tistory.com/293"

# Row 37
$ws.Range("D37").Value = "[Paper Review] DisenHAN: Disentangled Heterogeneous Graph Attention Network for Recommendation"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1436&mod=document&pageid=1"

# Row 39
$ws.Range("D39").Value = "Deep Face Detection with MTCNN in Python"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Deep-Face-Detection-with-MTCNN-in-Python-1"

# Row 43
$ws.Range("D43").Value = "신경레벨 from 정형외과 밀러책"
$ws.Range("E43").Value = "https://nittaku.tistory.com/506"

# Row 45
$ws.Range("D45").Value = "return_state, return_sequences"
$ws.Range("E45").Value = "https://dive-into-ds.tistory.com/80"

# Row 46
$ws.Range("D46").Value = "심장전도시스템 (cardiac conduction system)"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/371"
